$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.1657811205202583
$ws.Range("D2").Value = 0.1657810991108276

$ws.Range("C3").Value = 0.1615951055592954
$ws.Range("D3").Value = 0.1615950709247715

$ws.Range("C4").Value = 0.124265581389829
$ws.Range("D4").Value = 0.1242655813898289

$ws.Range("C5").Value = 0.124265581389829
$ws.Range("D5").Value = 0.124265581389829

$ws.Range("C6").Value = 0.124265581389829
$ws.Range("D6").Value = 0.124265581389829

$ws.Range("C7").Value = 0.1525293361577476
$ws.Range("D7").Value = 0.1525293532280509

$ws.Range("C8").Value = 0.1472976935932118
$ws.Range("D8").Value = 0.1472977325668632
